$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantity values on row 2
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 2

# Move the active selection to D3 (as reflected in the saved view state)
$ws.Range("D3").Select()
